$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (employee_id, employee_name, department, absence_reason, absence_duration, absence_date_serial, salary)
$rows = @(
    @{Row=2;  A=11834; B="Dra. Maria Eduarda Garcia"; C="Marketing";         D="Outros";              E=5; F=45084; G=2480.15}
    @{Row=3;  A=6591;  B="Lunna da Rosa";              C="Juridico";          D="Viagem de negocios";  E=8; F=45089; G=5984.7}
    @{Row=4;  A=63847; B="Evelyn Cirino";              C="Recursos Humanos";  D="Doenca";              E=4; F=45088; G=2874.31}
    @{Row=5;  A=24142; B="Dra. Nina Barros";           C="TI";                D="Outros";              E=8; F=45082; G=7044.78}
    @{Row=6;  A=46104; B="Cauã Moraes";                C="TI";                D="Viagem de negocios";  E=7; F=45103; G=5166.09}
    @{Row=7;  A=27727; B="Lucas Gabriel Alves";        C="Engenharia";        D="Viagem de negocios";  E=3; F=45104; G=5065.64}
    @{Row=8;  A=3517;  B="Srta. Sophie Nascimento";    C="P&D";               D="Viagem de negocios";  E=7; F=45093; G=6395.67}
    @{Row=9;  A=31830; B="Vitor Gabriel Duarte";       C="Engenharia";        D="Problemas pessoais";  E=5; F=45079; G=4478.25}
    @{Row=10; A=68160; B="Sr. João Felipe Cirino";     C="TI";                D="Consulta medica";     E=7; F=45097; G=6699.06}
    @{Row=11; A=12163; B="Davi Miguel Cirino";         C="Recursos Humanos";  D="Consulta medica";     E=7; F=45098; G=8772.09}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
